$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.228519
$ws.Range("H2").Value = 9.685557
$ws.Range("I2").Value = 0.0641296566303666
$ws.Range("J2").Value = 0.0641296566303666
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.33599166666667
$ws.Range("N2").Value = 52.007975
$ws.Range("O2").Value = 0.4573561888773979
$ws.Range("P2").Value = 0.4573561888773979
$ws.Range("Q2").Value = 55.96957847967501
$ws.Range("R2").Value = 503.726206317075
$ws.Range("S2").Value = 0.02933009535048062
$ws.Range("T2").Value = 0.02933009535048062

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.228519
$ws.Range("H3").Value = 9.685557
$ws.Range("I3").Value = 0.0641296566303666
$ws.Range("J3").Value = 0.0641296566303666
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.077707333333334
$ws.Range("N3").Value = 27.233122
$ws.Range("O3").Value = 0.2394870573052156
$ws.Range("P3").Value = 0.2394870573052156
$ws.Range("Q3").Value = 29.307550602106
$ws.Range("R3").Value = 263.767955418954
$ws.Range("S3").Value = 0.01535822275240041
$ws.Range("T3").Value = 0.0153582227524004

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.228519
$ws.Range("H4").Value = 9.685557
$ws.Range("I4").Value = 0.0641296566303666
$ws.Range("J4").Value = 0.0641296566303666
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.491094
$ws.Range("N4").Value = 34.473282
$ws.Range("O4").Value = 0.3031567538173866
$ws.Range("P4").Value = 0.3031567538173866
$ws.Range("Q4").Value = 37.099215309786
$ws.Range("R4").Value = 333.8929377880739
$ws.Range("S4").Value = 0.01944133852748558
$ws.Range("T4").Value = 0.01944133852748558

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.281951000000001
$ws.Range("H5").Value = 27.845853
$ws.Range("I5").Value = 0.1843719459262553
$ws.Range("J5").Value = 0.1843719459262553
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.33599166666667
$ws.Range("N5").Value = 52.007975
$ws.Range("O5").Value = 0.4573561888773979
$ws.Range("P5").Value = 0.4573561888773979
$ws.Range("Q5").Value = 160.9118251864084
$ws.Range("R5").Value = 1448.206426677675
$ws.Range("S5").Value = 0.08432365052474182
$ws.Range("T5").Value = 0.08432365052474182

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 9.281951000000001
$ws.Range("H6").Value = 27.845853
$ws.Range("I6").Value = 0.1843719459262553
$ws.Range("J6").Value = 0.1843719459262553
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.077707333333334
$ws.Range("N6").Value = 27.233122
$ws.Range("O6").Value = 0.2394870573052156
$ws.Range("P6").Value = 0.2394870573052156
$ws.Range("Q6").Value = 84.25883466034068
$ws.Range("R6").Value = 758.3295119430661
$ws.Range("S6").Value = 0.04415469477951522
$ws.Range("T6").Value = 0.04415469477951522

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9.281951000000001
$ws.Range("H7").Value = 27.845853
$ws.Range("I7").Value = 0.1843719459262553
$ws.Range("J7").Value = 0.1843719459262553
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.491094
$ws.Range("N7").Value = 34.473282
$ws.Range("O7").Value = 0.3031567538173866
$ws.Range("P7").Value = 0.3031567538173866
$ws.Range("Q7").Value = 106.659771444394
$ws.Range("R7").Value = 959.937942999546
$ws.Range("S7").Value = 0.0558936006219983
$ws.Range("T7").Value = 0.0558936006219983

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 37.833149
$ws.Range("H8").Value = 113.499447
$ws.Range("I8").Value = 0.751498397443378
$ws.Range("J8").Value = 0.7514983974433781
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 17.33599166666667
$ws.Range("N8").Value = 52.007975
$ws.Range("O8").Value = 0.4573561888773979
$ws.Range("P8").Value = 0.4573561888773979
$ws.Range("Q8").Value = 655.8751557877583
$ws.Range("R8").Value = 5902.876402089825
$ws.Range("S8").Value = 0.3437024430021754
$ws.Range("T8").Value = 0.3437024430021755

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 37.833149
$ws.Range("H9").Value = 113.499447
$ws.Range("I9").Value = 0.751498397443378
$ws.Range("J9").Value = 0.7514983974433781
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.077707333333334
$ws.Range("N9").Value = 27.233122
$ws.Range("O9").Value = 0.2394870573052156
$ws.Range("P9").Value = 0.2394870573052156
$ws.Range("Q9").Value = 343.4382541203927
$ws.Range("R9").Value = 3090.944287083534
$ws.Range("S9").Value = 0.1799741397733
$ws.Range("T9").Value = 0.1799741397733

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 37.833149
$ws.Range("H10").Value = 113.499447
$ws.Range("I10").Value = 0.751498397443378
$ws.Range("J10").Value = 0.7514983974433781
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 11.491094
$ws.Range("N10").Value = 34.473282
$ws.Range("O10").Value = 0.3031567538173866
$ws.Range("P10").Value = 0.3031567538173866
$ws.Range("Q10").Value = 434.7442714750059
$ws.Range("R10").Value = 3912.698443275054
$ws.Range("S10").Value = 0.2278218146679027
$ws.Range("T10").Value = 0.2278218146679027
